# Fix column order in the pf_bus_sc_results_all_cases.xlsx sheets.
#
# The sheets that carry zero-sequence / positive-sequence / negative-sequence
# short circuit impedance results (LL_*, LLG_*, LG_*) had their rk/xk columns
# grouped as: rk0, rk1, rk2, xk0, xk1, xk2 (columns H..M).
#
# The corrected layout interleaves them by sequence number instead:
#           rk0, xk0, rk1, xk1, rk2, xk2   (columns H..M)
#
# Column H (rk0) and M (xk2) are unaffected; columns I, J, K, L need their
# header labels and their 4 data rows (2-5) rearranged:
#   new I (pf_xk0_ohm) <- old K (pf_xk0_ohm data)
#   new J (pf_rk1_ohm) <- old I (pf_rk1_ohm data)
#   new K (pf_xk1_ohm) <- old L (pf_xk1_ohm data)
#   new L (pf_rk2_ohm) <- old J (pf_rk2_ohm data)

$wb = $excel.ActiveWorkbook

$affectedSheets = @(
    "LL_max", "LL_max_fault", "LL_min", "LL_min_fault",
    "LLG_max", "LLG_max_fault", "LLG_min", "LLG_min_fault",
    "LG_max", "LG_max_fault", "LG_min", "LG_min_fault"
)

foreach ($sheetName in $affectedSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Fix header row (row 1): I/J/K/L text labels.
    $ws.Range("I1").Value = "pf_xk0_ohm"
    $ws.Range("J1").Value = "pf_rk1_ohm"
    $ws.Range("K1").Value = "pf_xk1_ohm"
    $ws.Range("L1").Value = "pf_rk2_ohm"

    # Fix data rows 2-5: rearrange I/J/K/L values.
    for ($row = 2; $row -le 5; $row++) {
        $oldI = $ws.Range("I$row").Value()
        $oldJ = $ws.Range("J$row").Value()
        $oldK = $ws.Range("K$row").Value()
        $oldL = $ws.Range("L$row").Value()

        $ws.Range("I$row").Value = $oldK
        $ws.Range("J$row").Value = $oldI
        $ws.Range("K$row").Value = $oldL
        $ws.Range("L$row").Value = $oldJ
    }
}
